$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the existing header cell (H1) onto the new header cells,
# then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value2 = "I0"
$ws.Range("J1").Value2 = "IF"

# Fill in the new I (always 1) and J (mirrors H) data columns for rows 2-29.
for ($r = 2; $r -le 29; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value2 = 1
    $ws.Cells.Item($r, 10).Value2 = $hVal
}
